$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.9394706710291985
$ws.Cells.Item(2, 4).Value = 0.006976860367000626
$ws.Cells.Item(2, 5).Value = 0.0378358063093529
$ws.Cells.Item(2, 6).Value = 3.517797489769237
$ws.Cells.Item(2, 7).Value = 0.002607396914752075
$ws.Cells.Item(2, 9).Value = 1.821690953948121
$ws.Cells.Item(2, 10).Value = 0.08682373588464465
$ws.Cells.Item(2, 11).Value = 1.775038974064472
$ws.Cells.Item(2, 12).Value = 0.6367092987006515
$ws.Cells.Item(2, 14).Value = 3.231951835829115
$ws.Cells.Item(3, 2).Value = 0.9200326440561355
$ws.Cells.Item(3, 4).Value = 0.006603349257552793
$ws.Cells.Item(3, 5).Value = 0.03733750112965506
$ws.Cells.Item(3, 6).Value = 3.510760016306861
$ws.Cells.Item(3, 7).Value = 0.002612074777057436
$ws.Cells.Item(3, 9).Value = 1.826733587593125
$ws.Cells.Item(3, 10).Value = 0.08624461744051892
$ws.Cells.Item(3, 11).Value = 1.70483193987522
$ws.Cells.Item(3, 12).Value = 0.6248057736015653
$ws.Cells.Item(3, 14).Value = 3.251186447772795
$ws.Cells.Item(4, 2).Value = 0.9085689042117906
$ws.Cells.Item(4, 4).Value = 0.006371484654035697
$ws.Cells.Item(4, 5).Value = 0.03702554912446132
$ws.Cells.Item(4, 6).Value = 3.507965252309731
$ws.Cells.Item(4, 7).Value = 0.002615100406385217
$ws.Cells.Item(4, 9).Value = 1.830456976551218
$ws.Cells.Item(4, 10).Value = 0.08588174631521461
$ws.Cells.Item(4, 11).Value = 1.662768302081815
$ws.Cells.Item(4, 12).Value = 0.6178391480011101
$ws.Cells.Item(4, 14).Value = 3.26380851133078
$ws.Cells.Item(5, 2).Value = 0.9040161694374262
$ws.Cells.Item(5, 4).Value = 0.006276326617999928
$ws.Cells.Item(5, 5).Value = 0.03689689758373227
$ws.Cells.Item(5, 6).Value = 3.50721000263772
$ws.Cells.Item(5, 7).Value = 0.002616372074329108
$ws.Cells.Item(5, 9).Value = 1.832132023001506
$ws.Cells.Item(5, 10).Value = 0.08573202252093459
$ws.Cells.Item(5, 11).Value = 1.645889178881333
$ws.Cells.Item(5, 12).Value = 0.6150862048524317
$ws.Cells.Item(5, 14).Value = 3.269156191708689
$ws.Cells.Item(6, 2).Value = 0.9032673780700122
$ws.Cells.Item(6, 4).Value = 0.006260484110869768
$ws.Cells.Item(6, 5).Value = 0.03687544212638194
$ws.Cells.Item(6, 6).Value = 3.507107760567891
$ws.Cells.Item(6, 7).Value = 0.002616585574795886
$ws.Cells.Item(6, 9).Value = 1.832419690735087
$ws.Cells.Item(6, 10).Value = 0.08570704869522316
$ws.Cells.Item(6, 11).Value = 1.643102237806744
$ws.Cells.Item(6, 12).Value = 0.6146342748753
$ws.Cells.Item(6, 14).Value = 3.270056494027429
$ws.Cells.Item(7, 2).Value = 0.9085070228616416
$ws.Cells.Item(7, 4).Value = 0.006370204084134912
$ws.Cells.Item(7, 5).Value = 0.03702382030266271
$ws.Cells.Item(7, 6).Value = 3.507953513624926
$ws.Cells.Item(7, 7).Value = 0.002615117399647361
$ws.Cells.Item(7, 9).Value = 1.830478928126638
$ws.Cells.Item(7, 10).Value = 0.0858797346046174
$ws.Cells.Item(7, 11).Value = 1.662539603294618
$ws.Cells.Item(7, 12).Value = 0.6178016726239406
$ws.Cells.Item(7, 14).Value = 3.263879805800649
$ws.Cells.Item(8, 2).Value = 0.9326708081197239
$ws.Cells.Item(8, 4).Value = 0.006848576549902674
$ws.Cells.Item(8, 5).Value = 0.03766522055178623
$ws.Cells.Item(8, 6).Value = 3.515054085556315
$ws.Cells.Item(8, 7).Value = 0.002608978078115659
$ws.Cells.Item(8, 9).Value = 1.823299491558728
$ws.Cells.Item(8, 10).Value = 0.08662555808898631
$ws.Cells.Item(8, 11).Value = 1.750614737520891
$ws.Cells.Item(8, 12).Value = 0.6325339166372146
$ws.Cells.Item(8, 14).Value = 3.238415365038492
$ws.Cells.Item(9, 2).Value = 0.9837848492852572
$ws.Cells.Item(9, 4).Value = 0.007768261333275461
$ws.Cells.Item(9, 5).Value = 0.03887651204996345
$ws.Cells.Item(9, 6).Value = 3.541100543492462
$ws.Cells.Item(9, 7).Value = 0.002598150368226192
$ws.Cells.Item(9, 9).Value = 1.81419717539243
$ws.Cells.Item(9, 10).Value = 0.08803111890760817
$ws.Cells.Item(9, 11).Value = 1.931637466542156
$ws.Cells.Item(9, 12).Value = 0.6641428359853592
$ws.Cells.Item(9, 14).Value = 3.194925418027012
$ws.Cells.Item(10, 2).Value = 1.023602895980247
$ws.Cells.Item(10, 4).Value = 0.008435070263981714
$ws.Cells.Item(10, 5).Value = 0.03973960263871223
$ws.Cells.Item(10, 6).Value = 3.567650733263918
$ws.Cells.Item(10, 7).Value = 0.002590925815742676
$ws.Cells.Item(10, 9).Value = 1.81054518348575
$ws.Cells.Item(10, 10).Value = 0.08903027771791905
$ws.Cells.Item(10, 11).Value = 2.069752200506912
$ws.Cells.Item(10, 12).Value = 0.6890325281724188
$ws.Cells.Item(10, 14).Value = 3.166907277455152
$ws.Cells.Item(11, 2).Value = 1.042207446416
$ws.Cells.Item(11, 4).Value = 0.008737002172452435
$ws.Cells.Item(11, 5).Value = 0.0401267441132136
$ws.Cells.Item(11, 6).Value = 3.581344705936374
$ws.Cells.Item(11, 7).Value = 0.002587796113480404
$ws.Cells.Item(11, 9).Value = 1.809543443138857
$ws.Cells.Item(11, 10).Value = 0.08947781781954234
$ws.Cells.Item(11, 11).Value = 2.133707936777512
$ws.Cells.Item(11, 12).Value = 0.7007196015642592
$ws.Cells.Item(11, 14).Value = 3.155016361132809
$ws.Cells.Item(12, 2).Value = 1.049322899970122
$ws.Cells.Item(12, 4).Value = 0.008851175693916247
$ws.Cells.Item(12, 5).Value = 0.0402725817599725
$ws.Cells.Item(12, 6).Value = 3.586763014009904
$ws.Cells.Item(12, 7).Value = 0.002586633392486748
$ws.Cells.Item(12, 9).Value = 1.809258981875431
$ws.Cells.Item(12, 10).Value = 0.08964630627126802
$ws.Cells.Item(12, 11).Value = 2.158089005847444
$ws.Cells.Item(12, 12).Value = 0.7051977386696819
$ws.Cells.Item(12, 14).Value = 3.150636603373741
$ws.Cells.Item(13, 2).Value = 1.04778733737021
$ws.Cells.Item(13, 4).Value = 0.008826592922446252
$ws.Cells.Item(13, 5).Value = 0.04024120662715447
$ws.Cells.Item(13, 6).Value = 3.58558573029427
$ws.Cells.Item(13, 7).Value = 0.002586882809547949
$ws.Cells.Item(13, 9).Value = 1.809316025534393
$ws.Cells.Item(13, 10).Value = 0.08961006278819816
$ws.Cells.Item(13, 11).Value = 2.152830876073494
$ws.Cells.Item(13, 12).Value = 0.704230955901977
$ws.Cells.Item(13, 14).Value = 3.151574386002551
$ws.Cells.Item(14, 2).Value = 1.042791431617275
$ws.Cells.Item(14, 4).Value = 0.008746398302488245
$ws.Cells.Item(14, 5).Value = 0.04013875744627793
$ws.Cells.Item(14, 6).Value = 3.58178580806053
$ws.Cells.Item(14, 7).Value = 0.002587700006817593
$ws.Cells.Item(14, 9).Value = 1.809518138695111
$ws.Cells.Item(14, 10).Value = 0.08949169910355614
$ws.Cells.Item(14, 11).Value = 2.135710524480373
$ws.Cells.Item(14, 12).Value = 0.7010869681267025
$ws.Cells.Item(14, 14).Value = 3.154653568229179
$ws.Cells.Item(15, 2).Value = 1.039740442964586
$ws.Cells.Item(15, 4).Value = 0.008697256945922049
$ws.Cells.Item(15, 5).Value = 0.04007590552990514
$ws.Cells.Item(15, 6).Value = 3.579488559055719
$ws.Cells.Item(15, 7).Value = 0.002588203481347312
$ws.Cells.Item(15, 9).Value = 1.809654295471262
$ws.Cells.Item(15, 10).Value = 0.08941907017786654
$ws.Cells.Item(15, 11).Value = 2.12524497899841
$ws.Cells.Item(15, 12).Value = 0.6991680232996487
$ws.Cells.Item(15, 14).Value = 3.156555689357447
$ws.Cells.Item(16, 2).Value = 1.022396901007312
$ws.Cells.Item(16, 4).Value = 0.008415313497341259
$ws.Cells.Item(16, 5).Value = 0.03971419375069196
$ws.Cells.Item(16, 6).Value = 3.566788346276439
$ws.Cells.Item(16, 7).Value = 0.002591133494002627
$ws.Cells.Item(16, 9).Value = 1.810623924356591
$ws.Cells.Item(16, 10).Value = 0.08900089101940623
$ws.Cells.Item(16, 11).Value = 2.065595248462785
$ws.Cells.Item(16, 12).Value = 0.6882760929906624
$ws.Cells.Item(16, 14).Value = 3.167701589724658
$ws.Cells.Item(17, 2).Value = 1.011882766958109
$ws.Cells.Item(17, 4).Value = 0.008242019457014749
$ws.Cells.Item(17, 5).Value = 0.03949090801441368
$ws.Cells.Item(17, 6).Value = 3.559411321605594
$ws.Cells.Item(17, 7).Value = 0.00259297103520754
$ws.Cells.Item(17, 9).Value = 1.811387712714421
$ws.Cells.Item(17, 10).Value = 0.08874257662825258
$ws.Cells.Item(17, 11).Value = 2.029290919751588
$ws.Cells.Item(17, 12).Value = 0.6816876801707679
$ws.Cells.Item(17, 14).Value = 3.174758262966023
$ws.Cells.Item(18, 2).Value = 1.005881560019475
$ws.Cells.Item(18, 4).Value = 0.00814221083410871
$ws.Cells.Item(18, 5).Value = 0.03936196245180223
$ws.Cells.Item(18, 6).Value = 3.55532035430781
$ws.Cells.Item(18, 7).Value = 0.002594042704051565
$ws.Cells.Item(18, 9).Value = 1.811889103330081
$ws.Cells.Item(18, 10).Value = 0.08859334268683483
$ws.Cells.Item(18, 11).Value = 2.008515650574338
$ws.Cells.Item(18, 12).Value = 0.6779325158712766
$ws.Cells.Item(18, 14).Value = 3.178897510331765
$ws.Cells.Item(19, 2).Value = 1.003857606527163
$ws.Cells.Item(19, 4).Value = 0.008108393037858264
$ws.Cells.Item(19, 5).Value = 0.0393182140899091
$ws.Cells.Item(19, 6).Value = 3.553961337884488
$ws.Cells.Item(19, 7).Value = 0.002594408092116973
$ws.Cells.Item(19, 9).Value = 1.812069526935559
$ws.Cells.Item(19, 10).Value = 0.08854270100002104
$ws.Cells.Item(19, 11).Value = 2.001499698278622
$ws.Cells.Item(19, 12).Value = 0.6766669737862969
$ws.Cells.Item(19, 14).Value = 3.180312798058253
$ws.Cells.Item(20, 2).Value = 1.012997230561837
$ws.Cells.Item(20, 4).Value = 0.008260480530477565
$ws.Cells.Item(20, 5).Value = 0.03951473054891164
$ws.Cells.Item(20, 6).Value = 3.56018087440961
$ws.Cells.Item(20, 7).Value = 0.002592773898473495
$ws.Cells.Item(20, 9).Value = 1.811299980929071
$ws.Cells.Item(20, 10).Value = 0.08877014263484284
$ws.Cells.Item(20, 11).Value = 2.03314460189614
$ws.Cells.Item(20, 12).Value = 0.6823854758948471
$ws.Cells.Item(20, 14).Value = 3.173998742082446
$ws.Cells.Item(21, 2).Value = 1.044256944376087
$ws.Cells.Item(21, 4).Value = 0.008769957467290368
$ws.Cells.Item(21, 5).Value = 0.04016886980810419
$ws.Cells.Item(21, 6).Value = 3.582895619569001
$ws.Cells.Item(21, 7).Value = 0.002587459368254275
$ws.Cells.Item(21, 9).Value = 1.809456197990301
$ws.Cells.Item(21, 10).Value = 0.0895264919765193
$ws.Cells.Item(21, 11).Value = 2.14073477539614
$ws.Cells.Item(21, 12).Value = 0.7020090082999388
$ws.Cells.Item(21, 14).Value = 3.1537457961672
$ws.Cells.Item(22, 2).Value = 1.065096678179174
$ws.Cells.Item(22, 4).Value = 0.009102007536206713
$ws.Cells.Item(22, 5).Value = 0.04059194326482363
$ws.Cells.Item(22, 6).Value = 3.59909750418592
$ws.Cells.Item(22, 7).Value = 0.002584116696958536
$ws.Cells.Item(22, 9).Value = 1.808804206025776
$ws.Cells.Item(22, 10).Value = 0.0900150759109799
$ws.Cells.Item(22, 11).Value = 2.211998079421448
$ws.Cells.Item(22, 12).Value = 0.7151401400555528
$ws.Cells.Item(22, 14).Value = 3.141226779628042
$ws.Cells.Item(23, 2).Value = 1.053936733081542
$ws.Cells.Item(23, 4).Value = 0.008924857167603051
$ws.Cells.Item(23, 5).Value = 0.04036653966962866
$ws.Cells.Item(23, 6).Value = 3.590326034432763
$ws.Cells.Item(23, 7).Value = 0.00258588882438946
$ws.Cells.Item(23, 9).Value = 1.809101573439875
$ws.Cells.Item(23, 10).Value = 0.08975482789327849
$ws.Cells.Item(23, 11).Value = 2.17387673233327
$ws.Cells.Item(23, 12).Value = 0.7081037854121632
$ws.Cells.Item(23, 14).Value = 3.147842708791543
$ws.Cells.Item(24, 2).Value = 1.012493246019744
$ws.Cells.Item(24, 4).Value = 0.008252134842635428
$ws.Cells.Item(24, 5).Value = 0.03950396217642371
$ws.Cells.Item(24, 6).Value = 3.559832491801188
$ws.Cells.Item(24, 7).Value = 0.002592862976464721
$ws.Cells.Item(24, 9).Value = 1.811339450476915
$ws.Cells.Item(24, 10).Value = 0.08875768230514502
$ws.Cells.Item(24, 11).Value = 2.031402051419718
$ws.Cells.Item(24, 12).Value = 0.6820699009028033
$ws.Cells.Item(24, 14).Value = 3.174341865208106
$ws.Cells.Item(25, 2).Value = 0.9695588521705929
$ws.Cells.Item(25, 4).Value = 0.007521192775790553
$ws.Cells.Item(25, 5).Value = 0.03855364894939761
$ws.Cells.Item(25, 6).Value = 3.532754090591681
$ws.Cells.Item(25, 7).Value = 0.002600950681446724
$ws.Cells.Item(25, 9).Value = 1.81612666107965
$ws.Cells.Item(25, 10).Value = 0.08765685175759863
$ws.Cells.Item(25, 11).Value = 1.881771223972862
$ws.Cells.Item(25, 12).Value = 0.6552997742508069
$ws.Cells.Item(25, 14).Value = 3.206000101529995
